$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.538.53"
$ws.Range("E2").Value = "  +6.88%  "

$ws.Range("D3").Value = "1.725.37"
$ws.Range("E3").Value = "  +3.67%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "332.83"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").Value = "0.3718"
$ws.Range("E7").Value = "  +1.95%  "

$ws.Range("D8").Value = "48.46"
$ws.Range("E8").Value = "  +2.59%  "

$ws.Range("D9").Value = "0.3365"
$ws.Range("E9").Value = "  +3.40%  "

$ws.Range("E10").Value = "  +4.27%  "

$ws.Range("D11").Value = "0.07417"
$ws.Range("E11").Value = "  +4.96%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").Value = "6.385"
$ws.Range("E13").Value = "  +5.29%  "

$ws.Range("D14").Value = "20.10"
$ws.Range("E14").Value = "  +3.21%  "

$ws.Range("D15").Value = "7.043"
$ws.Range("E15").Value = "  +6.90%  "

$ws.Range("D16").Value = "1.726.64"
$ws.Range("E16").Value = "  +3.58%  "

$ws.Range("D17").Value = "0.00001072"
$ws.Range("E17").Value = "  +2.40%  "

$ws.Range("D18").Value = "0.06637"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "81.85"
$ws.Range("E19").Value = "  +4.34%  "

$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").Value = "16.52"

$ws.Range("D22").Value = "6.142"
$ws.Range("E22").Value = "  +3.75%  "

$ws.Range("D23").Value = "12.77"
$ws.Range("E23").Value = "  +2.14%  "

$ws.Range("D24").Value = "26.528.24"
$ws.Range("E24").Value = "  +6.95%  "

$ws.Range("D25").Value = "2.453"
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("D26").Value = "1.411"
$ws.Range("E26").Value = "  +21.06%  "

$ws.Range("D27").Value = "2.394"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").Value = "150.88"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").Value = "19.39"
$ws.Range("E29").Value = "  +4.07%  "

$ws.Range("D30").Value = "1.920.51"
$ws.Range("E30").Value = "  +3.82%  "

$ws.Range("D31").Value = "131.14"

$ws.Range("D32").Value = "4.111"
$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("D33").Value = "5.964"
$ws.Range("E33").Value = "  +4.92%  "

$ws.Range("D34").Value = "0.08626"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("D35").Value = "1.695"
$ws.Range("E35").Value = "  +2.94%  "

$ws.Range("E36").Value = "  +4.90%  "

$ws.Range("E37").Value = "  +3.80%  "

$ws.Range("E38").Value = "  +2.22%  "

$ws.Range("D39").Value = "0.06207"
$ws.Range("E39").Value = "  -0.69%  "

$ws.Range("D40").Value = "0.2151"
$ws.Range("E40").Value = "  +3.19%  "

$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").Value = "1.224"
$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("D43").Value = "0.6205"
$ws.Range("E43").Value = "  +4.83%  "

$ws.Range("D44").Value = "14.14"
$ws.Range("E44").Value = "  +6.31%  "

$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").Value = "3.835"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "0.6014"
$ws.Range("E47").Value = "  +6.33%  "

$ws.Range("D48").Value = "128.88"
$ws.Range("E48").Value = "  +2.56%  "

$ws.Range("D49").Value = "2.042"
$ws.Range("E49").Value = "  +4.97%  "

$ws.Range("D50").Value = "0.07169"
$ws.Range("E50").Value = "  +2.79%  "

$ws.Range("D51").Value = "77.01"
$ws.Range("E51").Value = "  +2.58%  "
